$d = $word.ActiveDocument

# The document body already contains runs that reference the built-in
# character styles "Strong" (bold) and "Emphasis" (italic) via rStyle,
# but the corresponding <w:style> definitions are missing from
# styles.xml. Materialize both style definitions (universal/built-in
# style mapping) so the references resolve correctly.

# --- Strong (bold) ---
$strong = $d.Styles.Add("Strong", 2)              # wdStyleTypeCharacter
$strong.BaseStyle = "Fuentedeprrafopredeter"       # Default Paragraph Font
$strong.Priority = 22
$strong.QuickStyle = $true
$strong.Font.Bold = $true
$strong.Font.BoldBi = $true

# --- Emphasis (italic) ---
$emphasis = $d.Styles.Add("Emphasis", 2)           # wdStyleTypeCharacter
$emphasis.BaseStyle = "Fuentedeprrafopredeter"     # Default Paragraph Font
$emphasis.Priority = 20
$emphasis.QuickStyle = $true
$emphasis.Font.Italic = $true
$emphasis.Font.ItalicBi = $true

Write-Host "Added Strong and Emphasis character styles"
